$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing rows 2..45 down to 3..46)
$ws.Rows.Item(2).Insert()

# Copy the (now shifted) row 3 formatting into the new blank row 2 so the
# inserted row matches the table's look (centered data cells, price number
# format, etc.) instead of the bold header formatting Excel defaults to.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Fill in the new latest-price entry
$ws.Range("A2").Value = 45
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 285.25
$ws.Range("E2").Value = "24.10.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-24-october-2025.pdf"

# Row 26 (previously row 25 before the insert) gains the hyperlink that the
# source row never had wired up in the workbook.
$ws.Hyperlinks.Add($ws.Range("F26"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Restore the plain data-cell formatting on F26 (Hyperlinks.Add stamps the
# built-in blue/underline Hyperlink style by default).
$ws.Range("E26").Copy()
$ws.Range("F26").PasteSpecial(-4122)

Write-Output "edit complete"
